# Append a new "Job Applications" time-tracking entry as row 3
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = " 10/11/2025"
$ws.Range("B3").Value = "16:46:55"
$ws.Range("C3").Value = "00:00:07"
$ws.Range("D3").Value = "Job Applications`n"

# Keep the row height at its default (matches rows above, which also
# contain embedded newlines but no custom row height).
$ws.Rows(3).AutoFit()
